# Add the 24-October-2020 "mystery case" row to the VIC_Mystery_cases table.
# This is a new row inserted above the existing data (row 2), which pushes
# every other row down by one and grows the table/autofilter by a row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a blank row at the top of the data (row 2, just below the header)
# and give it the same look as the row it displaces by copying that row's
# formatting down into it.
$ws.Rows(2).Insert()
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)  # xlPasteFormats

# Grow the table definition (and its AutoFilter) to cover the extra row.
$lo.Resize($ws.Range("A1:E23"))

# Fill in the new record.
$ws.Range("A2").Value = 44125
$ws.Range("B2").Value = 3029
$ws.Range("C2").Value = "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-24-october-2020"
$ws.Range("D2").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]"
$ws.Range("E2").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]+14"

# The row that just got shifted past the table's old last row needs its
# calculated-column formulas re-asserted (structured references on that row
# otherwise come out malformed after the insert+resize).
$ws.Range("D23").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]"
$ws.Range("E23").Formula = "=VIC_Mystery_cases[[#This Row],[Date]]+14"

# The existing hyperlink (News Link column) belongs to the row that is now
# one below where it used to be; move it down to keep it on the same data.
$ws.Range("C4").Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-17-october-2020")
# Re-adding the hyperlink stamps the default hyperlink look on the cell;
# restore the plain table formatting that the other data cells use.
$ws.Range("C6").Copy()
$ws.Range("C5").PasteSpecial(-4122)  # xlPasteFormats

# Match the saved selection position.
$ws.Range("B2").Select() | Out-Null

Write-Host "Inserted new mystery case row (2020-10-24, postcode 3029)"
